$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: change Order type letters and quantities
$ws.Range("B2").Value = "w"
$ws.Range("C2").Value = "w"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("H2").Value = 18.2325

# Add new row 3 with an order number in column A
$ws.Range("A3").Value = 3
